# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# that used to follow the LOB1018 requirement line, collapsing the blank
# paragraph right after it with the blank paragraph that used to sit just
# before the page-break paragraph at the end of the document.

$d = $word.ActiveDocument

# Locate the LOB1018 requirement paragraph via Find (robust to position drift).
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "LOB1018: Física I (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the LOB1018 requirement paragraph."
}

$lob1018Index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $anchor.Start -and $p.Range.End -ge $anchor.End) {
        $lob1018Index = $i
        break
    }
}
if ($lob1018Index -eq -1) {
    throw "Could not resolve the paragraph index of the LOB1018 line."
}

$lob1018Para = $d.Paragraphs.Item($lob1018Index)

# The three paragraphs immediately following LOB1018 are:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: ... Creative Commons Attribution"
$blank = $lob1018Para.Next()
$jupiter = $blank.Next()
$copyright = $jupiter.Next()

if ($jupiter.Range.Text -notlike "Ver no Jupiter*") {
    throw "Unexpected paragraph where 'Ver no Jupiter...' was expected."
}
if ($copyright.Range.Text -notlike "*Contact: luizeleno@usp.br*") {
    throw "Unexpected paragraph where the copyright line was expected."
}

$deleteRange = $d.Range($blank.Range.Start, $copyright.Range.End)
$deleteRange.Delete()
